# Updates cryptos price/volume data to match the latest scrape.
# Every Price (D) / Volume (E) / Coin (B) / Link (C) cell in the source
# sheet is stored as text (t="inlineStr"), so Price cells are forced to
# Text format before assignment -- this keeps values such as "7.20" or
# "598.10" from being silently renumbered to 7.2 / 598.1 by Excel's
# automatic type inference, matching the original workbook's formatting.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "67.581.52"
$ws.Range("E2").Value = "  -2.30%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.808.58"
$ws.Range("E3").Value = "  +1.66%  "

$ws.Range("E4").Value = "  +0.08%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "597.98"
$ws.Range("E5").Value = "  -2.77%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "177.08"
$ws.Range("E6").Value = "  -0.59%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.807.39"
$ws.Range("E7").Value = "  +1.67%  "

$ws.Range("E8").Value = "  +0.03%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.526"
$ws.Range("E9").Value = "  -0.25%  "

$ws.Range("E10").Value = "  -4.40%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.19"
$ws.Range("E11").Value = "  -4.80%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.463"
$ws.Range("E12").Value = "  -4.25%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "38.36"
$ws.Range("E13").Value = "  -4.23%  "

$ws.Range("E14").Value = "  -3.51%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.441.84"
$ws.Range("E15").Value = "  +1.75%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.810.23"
$ws.Range("E16").Value = "  +1.63%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "67.661.45"
$ws.Range("E17").Value = "  -2.28%  "

$ws.Range("B18").Value = "Polkadot"
$ws.Range("C18").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.20"
$ws.Range("E18").Value = "  -3.51%  "

$ws.Range("B19").Value = "TRON"
$ws.Range("C19").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.115"
$ws.Range("E19").Value = "  -4.40%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "16.59"
$ws.Range("E20").Value = "  +1.33%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "491.40"
$ws.Range("E21").Value = "  -1.52%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.10"
$ws.Range("E22").Value = "  -2.83%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.742"
$ws.Range("E23").Value = "  +2.92%  "

$ws.Range("E24").Value = "  +13.84%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "84.92"
$ws.Range("E25").Value = "  -0.92%  "

$ws.Range("E26").Value = "  -6.60%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.31"
$ws.Range("E27").Value = "  -4.43%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.19"
$ws.Range("E28").Value = "  -5.91%  "

$ws.Range("E29").Value = "  +0.02%  "

$ws.Range("E30").Value = "  +0.92%  "

$ws.Range("E31").Value = "  -2.64%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "32.61"
$ws.Range("E32").Value = "  +6.97%  "

$ws.Range("E33").Value = "  -4.52%  "

$ws.Range("E34").Value = "  -4.10%  "

$ws.Range("E35").Value = "  +0.09%  "

$ws.Range("E36").Value = "  -3.69%  "

$ws.Range("E37").Value = "  -5.31%  "

$ws.Range("E38").Value = "  -2.45%  "

$ws.Range("E39").Value = "  -5.00%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "450.33"
$ws.Range("E40").Value = "  +0.05%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "49.19"
$ws.Range("E41").Value = "  -1.15%  "

$ws.Range("E42").Value = "  -2.97%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.89"
$ws.Range("E43").Value = "  -4.31%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.34"
$ws.Range("E44").Value = "  -2.77%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "41.23"
$ws.Range("E45").Value = "  -8.14%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.844.88"
$ws.Range("E46").Value = "  -3.47%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "139.61"
$ws.Range("E47").Value = "  +1.24%  "

$ws.Range("E49").Value = "  -2.65%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "24.72"
$ws.Range("E50").Value = "  +13.03%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "26.03"
$ws.Range("E51").Value = "  -5.12%  "
